$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder "Safety, Climate, Governance" -> "Safety, Governance, Climate" for rows 4-11
$rowsSCG = @(4,5,6,7,8,9,10,11)
foreach ($r in $rowsSCG) {
    $ws.Range("D$r").Value = "Safety, Governance, Climate"
}

# Reorder "Climate, Governance" -> "Governance, Climate" for rows 18,19,21,22,23,24,25
$rowsCG = @(18,19,21,22,23,24,25)
foreach ($r in $rowsCG) {
    $ws.Range("D$r").Value = "Governance, Climate"
}

# Remove the last data row (row 31) which held the "One Health Portal" RFP entry
$ws.Rows("31").Delete()

Write-Output "Edit applied"
